$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 (Nebraska) updates
$ws.Range("B28").Value = 44034
$ws.Range("C28").Value = 23486
$ws.Range("D28").Value = 311
$ws.Range("E28").Value = 1390
$ws.Range("G28").Value = 7.65
$ws.Range("H28").Value = 7.74
$ws.Range("K28").Value = 18180
$ws.Range("L28").Value = 297

# Row 41 (Iowa) updates
$ws.Range("C41").Value = 40085
$ws.Range("E41").Value = 3286

# Row 44 (Idaho) - previously empty/error, now filled in
$ws.Range("B44").Value = 44034
$ws.Range("B44").NumberFormat = "YYYY-MM-DD"
$ws.Range("C44").Value = 16322
$ws.Range("D44").Value = 131
$ws.Range("E44").Value = 208
$ws.Range("F44").Value = 2
$ws.Range("G44").Value = 1.27
$ws.Range("H44").Value = 1.53
$ws.Range("J44").Value = $true
$ws.Range("O44").Value = "Success!"
